$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 418 (a new week of price data,
# 2022-02-03 / serial 44595), pushing the existing rows 418-423 down to
# 421-426.
$ws.Rows.Item(418).Insert()
$ws.Rows.Item(418).Insert()
$ws.Rows.Item(418).Insert()

# Row 418 - Sandia, Extra
$ws.Cells.Item(418, 1).Value = 3
$ws.Cells.Item(418, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(418, 3).Value = "Coquimbo"
$ws.Cells.Item(418, 4).Value = 44595
$ws.Cells.Item(418, 5).Value = 5
$ws.Cells.Item(418, 6).Value = 100112028
$ws.Cells.Item(418, 7).Value = "Sandia"
$ws.Cells.Item(418, 8).Value = "Sin especificar"
$ws.Cells.Item(418, 9).Value = "Extra"
$ws.Cells.Item(418, 10).Value = 250
$ws.Cells.Item(418, 11).Value = 3000
$ws.Cells.Item(418, 12).Value = 3000
$ws.Cells.Item(418, 13).Value = 3000
$ws.Cells.Item(418, 14).Value = "$/unidad"
$ws.Cells.Item(418, 15).Value = "Paine"
$ws.Cells.Item(418, 16).Value = 3000
$ws.Cells.Item(418, 17).Value = 1
$ws.Cells.Item(418, 18).Value = "Hortaliza"

# Row 419 - Sandia, Primera
$ws.Cells.Item(419, 1).Value = 3
$ws.Cells.Item(419, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(419, 3).Value = "Coquimbo"
$ws.Cells.Item(419, 4).Value = 44595
$ws.Cells.Item(419, 5).Value = 5
$ws.Cells.Item(419, 6).Value = 100112028
$ws.Cells.Item(419, 7).Value = "Sandia"
$ws.Cells.Item(419, 8).Value = "Sin especificar"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 260
$ws.Cells.Item(419, 11).Value = 2000
$ws.Cells.Item(419, 12).Value = 2000
$ws.Cells.Item(419, 13).Value = 2000
$ws.Cells.Item(419, 14).Value = "$/unidad"
$ws.Cells.Item(419, 15).Value = "Paine"
$ws.Cells.Item(419, 16).Value = 2000
$ws.Cells.Item(419, 17).Value = 1
$ws.Cells.Item(419, 18).Value = "Hortaliza"

# Row 420 - Sandia, Segunda
$ws.Cells.Item(420, 1).Value = 3
$ws.Cells.Item(420, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(420, 3).Value = "Coquimbo"
$ws.Cells.Item(420, 4).Value = 44595
$ws.Cells.Item(420, 5).Value = 5
$ws.Cells.Item(420, 6).Value = 100112028
$ws.Cells.Item(420, 7).Value = "Sandia"
$ws.Cells.Item(420, 8).Value = "Sin especificar"
$ws.Cells.Item(420, 9).Value = "Segunda"
$ws.Cells.Item(420, 10).Value = 280
$ws.Cells.Item(420, 11).Value = 1500
$ws.Cells.Item(420, 12).Value = 1500
$ws.Cells.Item(420, 13).Value = 1500
$ws.Cells.Item(420, 14).Value = "$/unidad"
$ws.Cells.Item(420, 15).Value = "Paine"
$ws.Cells.Item(420, 16).Value = 1500
$ws.Cells.Item(420, 17).Value = 1
$ws.Cells.Item(420, 18).Value = "Hortaliza"
